$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 86, shifting existing rows 86:90 down to 87:91,
# then populate the new row 86 with the new data record.
$ws.Rows("86:86").Insert(-4121)

$ws.Range("A86").Value = 8
$ws.Range("B86").Value = "Terminal La Palmera de La Serena"
$ws.Range("C86").Value = "Coquimbo"
$ws.Range("D86").Value = 44615
$ws.Range("E86").Value = 4
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100109
$ws.Range("H86").Value = "Uva"
$ws.Range("I86").Value = 100109001
$ws.Range("J86").Value = "Uva"
$ws.Range("K86").Value = "Red Globe"
$ws.Range("L86").Value = "Primera"
$ws.Range("M86").Value = 400
$ws.Range("N86").Value = 11000
$ws.Range("O86").Value = 12000
$ws.Range("P86").Value = 11500
$ws.Range("Q86").Value = "`$/bandeja 18 kilos"
$ws.Range("R86").Value = "Provincia del Elquí"
$ws.Range("S86").Value = 639
$ws.Range("T86").Value = 18
